# VARILLA Y ALAMBRE VISILLO price list - Google Drive sync fix
# - bump the sheet date (A1) by one day
# - correct the VARILLA CHATA price (D22)
# - correct the ALAMBRE price (D38)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311

$ws.Range("D22").Value = 158.3

$ws.Range("D38").Value = 244.506
